$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Home")

# Insert two new rows above the existing "Author(s)" row to make room for
# the new "Policy Name" and "Description" metadata rows.
$ws.Rows("2:3").Insert()

# Copy the formatting from the (now shifted) "Author(s)" row down onto the
# two newly inserted rows so they pick up the same label style (s=31 on
# column A) instead of inheriting the title row's formatting.
$ws.Range("A4:B4").Copy()
$ws.Range("A2:B3").PasteSpecial(-4122)
$ws.Range("C2:D3").Clear()

# Populate the new metadata rows.
$ws.Range("A2").Value = "Policy Name"
$ws.Range("A3").Value = "Description"
$ws.Range("B3").Value = "HDescription"
$ws.Range("B2").Value = "HPolicyName"
